# Update logo, doc, and fix (#7)
#
# The deck's slide 3 ("SAIUnit" logo lockup) is reworked into the richer
# lockup (big "SAIUnit" wordmark + underline rule + "Unit-aware computation
# for Scientific AI" tagline) that slides 1 and 2 already use, and the
# original (plain) slide 3 is preserved as a new slide 4 right after it.

$p = $ppt.ActivePresentation

# Slides 1 and 2 already carry the target 4-shape lockup (textbox "holder",
# big wordmark textbox, tagline textbox, connector rule) at the exact
# target geometry, so duplicating slide 1 reproduces that structure (and
# its shape ids/names) exactly; it only needs its two text strings fixed
# up. The duplicate is inserted right after slide 1, so move it into slot
# 3 - this pushes the original (unmodified) slide 3 down to slot 4, which
# is exactly the "preserve the old slide" part of the edit.
$newSlide = $p.Slides.Item(1).Duplicate()
$newSlide.MoveTo(3)

# Fix the wordmark casing: "SaiUnit" -> "SAIUnit".
$wordmark = $newSlide.Shapes.Item(2)
$wordmark.TextFrame.TextRange.Text = "SAIUnit"

# Fix the tagline wording: "-aware Computations for " -> "-aware computation for ".
$tagline = $newSlide.Shapes.Item(3)
$full = $tagline.TextFrame.TextRange.Text
$oldFragment = "-aware Computations for "
$start = $full.IndexOf($oldFragment) + 1
$frag = $tagline.TextFrame.TextRange.Characters($start, $oldFragment.Length)
$frag.Text = "-aware computation for "
